$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) retains its text representation
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.308.81'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.930.03'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '249.77'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").Value = '0.7204'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '27.94'
$ws.Range("E8").Value = '  -2.98%  '
$ws.Range("D9").Value = '0.3204'
$ws.Range("E9").Value = '  -4.35%  '
$ws.Range("D10").Value = '0.07115'
$ws.Range("E10").Value = '  -3.72%  '
$ws.Range("D11").Value = '0.7882'
$ws.Range("E11").Value = '  -3.69%  '
$ws.Range("D12").Value = '0.08023'
$ws.Range("E12").Value = '  -1.55%  '
$ws.Range("D13").Value = '1.928.02'
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("D14").Value = '5.381'
$ws.Range("E14").Value = '  -2.11%  '
$ws.Range("D15").Value = '94.72'
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("E16").Value = '  -2.05%  '
$ws.Range("D17").Value = '30.302.79'
$ws.Range("E17").Value = '  -0.23%  '
$ws.Range("D18").Value = '256.60'
$ws.Range("E18").Value = '  +0.75%  '
$ws.Range("D19").Value = '0.000008080'
$ws.Range("E19").Value = '  -3.46%  '
$ws.Range("D20").Value = '5.739'
$ws.Range("E20").Value = '  -2.41%  '
$ws.Range("D21").Value = '2.183.74'
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("E24").Value = '  -2.65%  '
$ws.Range("D25").Value = '9.553'
$ws.Range("E25").Value = '  -3.52%  '
$ws.Range("D26").Value = '164.23'
$ws.Range("E26").Value = '  +1.71%  '
$ws.Range("E27").Value = '  -1.72%  '
$ws.Range("D28").Value = '2.301'
$ws.Range("E28").Value = '  -5.80%  '
$ws.Range("D29").Value = '0.1286'
$ws.Range("E29").Value = '  -2.24%  '
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("E31").Value = '  -2.54%  '
$ws.Range("E32").Value = '  -1.67%  '
$ws.Range("D33").Value = '4.155'
$ws.Range("E33").Value = '  -2.91%  '
$ws.Range("D34").Value = '0.05108'
$ws.Range("E34").Value = '  -3.93%  '
$ws.Range("D35").Value = '1.286'
$ws.Range("E35").Value = '  -2.76%  '
$ws.Range("D36").Value = '0.7501'
$ws.Range("E36").Value = '  -1.96%  '
$ws.Range("D37").Value = '2.769'
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("D38").Value = '0.01987'
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("D39").Value = '2.798'
$ws.Range("E39").Value = '  -1.75%  '
$ws.Range("D40").Value = '78.39'
$ws.Range("E40").Value = '  -3.82%  '
$ws.Range("D41").Value = '6.399'
$ws.Range("E41").Value = '  -2.92%  '
$ws.Range("D42").Value = '0.4523'
$ws.Range("E42").Value = '  -1.33%  '
$ws.Range("D43").Value = '1.995'
$ws.Range("E43").Value = '  -2.40%  '
$ws.Range("D44").Value = '0.8453'
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").Value = '0.9998'
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '101.41'
$ws.Range("E46").Value = '  -1.79%  '
$ws.Range("D47").Value = '9.828'
$ws.Range("E47").Value = '  -0.88%  '
$ws.Range("D48").Value = '7.491'
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").Value = '36.88'
$ws.Range("D50").Value = '964.62'
$ws.Range("D51").Value = '0.4212'
$ws.Range("E51").Value = '  -0.05%  '
